# Aula 11 - Leitura Ativa
# Commit: Corrige apresentacoes - logo em cada slide (Aulas 03-30)
# - embed the "Educa com Talento" logo individually on every slide
#   (previously it only lived on the now-unused MASTER slide layout)
# - give every slide its own dark-navy (1A1A2E) background fill
# - tighten up a number of bullet / title strings

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Text edits (content simplification) - done first, while shape indices
#    still match the pristine "before" ordering, and always addressed by
#    shape Name so ordering never matters anyway.
# ---------------------------------------------------------------------------

$p.Slides.Item(1).Shapes.Item("Text 3").TextFrame.TextRange.Text = "Técnicas de Leitura Ativa"

$p.Slides.Item(2).Shapes.Item("Text 5").TextFrame.TextRange.Text = "Aplicar técnica SQ3R"
$p.Slides.Item(2).Shapes.Item("Text 7").TextFrame.TextRange.Text = "Processar profundamente"

$p.Slides.Item(3).Shapes.Item("Text 5").TextFrame.TextRange.Text = "Mito de reler várias vezes"
$p.Slides.Item(3).Shapes.Item("Text 7").TextFrame.TextRange.Text = "Sublinhado não basta"

$p.Slides.Item(4).Shapes.Item("Text 7").TextFrame.TextRange.Text = "Read, Recite, Review"
$p.Slides.Item(4).Shapes.Item("Text 8").Delete()
$p.Slides.Item(4).Shapes.Item("Text 9").Delete()
$p.Slides.Item(4).Shapes.Item("Text 10").Delete()
$p.Slides.Item(4).Shapes.Item("Text 11").Delete()

$p.Slides.Item(5).Shapes.Item("Text 1").TextFrame.TextRange.Text = "Processamento"
$p.Slides.Item(5).Shapes.Item("Text 3").TextFrame.TextRange.Text = "Parafrasear"
$p.Slides.Item(5).Shapes.Item("Text 5").TextFrame.TextRange.Text = "Conectar conhecimento"

$p.Slides.Item(6).Shapes.Item("Text 1").TextFrame.TextRange.Text = "Sistema"
$p.Slides.Item(6).Shapes.Item("Text 7").TextFrame.TextRange.Text = "Pós: síntese"

$p.Slides.Item(7).Shapes.Item("Text 3").TextFrame.TextRange.Text = "Leitura ativa > passiva"
$p.Slides.Item(7).Shapes.Item("Text 7").TextFrame.TextRange.Text = "Processar: parafrasear, conectar"
$p.Slides.Item(7).Shapes.Item("Text 9").TextFrame.TextRange.Text = "Uma ativa = cinco passivas"

$p.Slides.Item(8).Shapes.Item("Text 3").TextFrame.TextRange.Text = "Aplicar SQ3R: texto 5-10 páginas + cada passo"

$p.Slides.Item(9).Shapes.Item("Text 3").TextFrame.TextRange.Text = '"Leitura ativa exige mais, mas economiza tempo."'

# ---------------------------------------------------------------------------
# 2) Per-slide background: solid fill 1A1A2E (was only set on the MASTER
#    slide layout before; now every slide carries its own override so the
#    deck still renders correctly once that layout goes away).
# ---------------------------------------------------------------------------

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.Background.Fill.Solid()
    $slide.Background.Fill.ForeColor.RGB = 0x2E1A1A
}

# ---------------------------------------------------------------------------
# 3) Logo picture on every slide: embed the "preencoded.png" logo (same
#    bytes the MASTER layout used to carry) at 274320,137160 EMU with size
#    1097280x457200 EMU -> 21.6,10.8 pt / 86.4,36.0 pt (1 pt = 12700 EMU),
#    then push it to the back of the z-order so it is the first shape in
#    the slide's shape tree, matching where the layout used to draw it.
# ---------------------------------------------------------------------------

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $logo = $slide.Shapes.AddPicture("educa_com_talento_logo.png", $false, $true, 21.6, 10.8, 86.4, 36.0)
    $logo.Name = "Image 0"
    $logo.ZOrder(1)
}
